$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 410; this shifts the existing rows 410-473
# down to 411-474 and extends the used range to A1:R474.
$ws.Rows.Item(410).Insert()

# Populate the newly inserted row 410 with the new weekly record.
$ws.Cells.Item(410, 1).Value = 2
$ws.Cells.Item(410, 2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(410, 3).Value = 'Coquimbo'
$ws.Cells.Item(410, 4).Value = '2023-10-19'
$ws.Cells.Item(410, 5).Value = 4
$ws.Cells.Item(410, 6).Value = 100112021
$ws.Cells.Item(410, 7).Value = 'Ají'
$ws.Cells.Item(410, 8).Value = 'Americana (o)'
$ws.Cells.Item(410, 9).Value = 'Primera'
$ws.Cells.Item(410, 10).Value = 160
$ws.Cells.Item(410, 11).Value = 35000
$ws.Cells.Item(410, 12).Value = 40000
$ws.Cells.Item(410, 13).Value = 37500
$ws.Cells.Item(410, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(410, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(410, 16).Value = 1500
$ws.Cells.Item(410, 17).Value = 25
$ws.Cells.Item(410, 18).Value = 'Hortaliza'
